$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the source inlineStr cells)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.483.54'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '1.572.86'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '291.54'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.3709'
$ws.Range("E7").Value = '  -1.58%  '
$ws.Range("D8").Value = '49.95'
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").Value = '0.07552'
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D13").Value = '21.27'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '6.039'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = '6.969'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = '1.571.07'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '0.00001122'
$ws.Range("E17").Value = '  -0.98%  '
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").Value = '0.06761'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '6.291'
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("E22").Value = '  -1.05%  '
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("D24").Value = '22.473.82'
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").Value = '2.355'
$ws.Range("E25").Value = '  -2.07%  '
$ws.Range("D26").Value = '2.631'
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").Value = '149.46'
$ws.Range("E28").Value = '  +1.23%  '
$ws.Range("D29").Value = '5.054'
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("D30").Value = '125.10'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = '1.746.78'
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").Value = '1.081'
$ws.Range("E32").Value = '  +9.39%  '
$ws.Range("D33").Value = '6.240'
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("D34").Value = '2.017'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '9.837'
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("D36").Value = '0.08358'
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").Value = '0.02488'
$ws.Range("E37").Value = '  -1.14%  '
$ws.Range("D38").Value = '0.2306'
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Value = '1.340'
$ws.Range("E39").Value = '  -4.70%  '
$ws.Range("D40").Value = '0.06559'
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").Value = '5.462'
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("D42").Value = '11.39'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '0.6227'
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '14.08'
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = '3.809'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("D48").Value = '129.08'
$ws.Range("E48").Value = '  +3.73%  '
$ws.Range("D49").Value = '2.074'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("E50").Value = '  -3.92%  '
$ws.Range("D51").Value = '0.07337'
$ws.Range("E51").Value = '  +0.35%  '
